$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain plain text,
# matching the inline-string cell type used in the original workbook,
# so strings that look numeric (e.g. "1.001") are not auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.569.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.653.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9971"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3645"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.57"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3258"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.128"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07037"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9962"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.972"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.607"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.654.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06608"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9972"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.78"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.943"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.75"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.60"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.588.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.459"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.345"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -16.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.61"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -9.00%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.229"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.837.32"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.52"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.066"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.737"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -16.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08454"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.666"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.24"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -11.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.276"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.211"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06029"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -9.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02230"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2069"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.131"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -11.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9970"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5902"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.813"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.60"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5627"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.949"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06946"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.192"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.99%  "
